$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 22, shifting existing rows 22-44 down to 23-45
$ws.Rows.Item(22).EntireRow.Insert()

# Populate the newly inserted row 22 with the new weekly data point
$ws.Range("A22").Value = 10
$ws.Range("B22").Value = "Vega Modelo de Temuco"
$ws.Range("C22").Value = "La Araucanía"
$ws.Range("D22").Value = 44447
$ws.Range("E22").Value = 9
$ws.Range("F22").Value = 100112035
$ws.Range("G22").Value = "Bruselas (repollito)"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 30
$ws.Range("K22").Value = 27000
$ws.Range("L22").Value = 27000
$ws.Range("M22").Value = 27000
$ws.Range("N22").Value = "$/malla 10 kilos"
$ws.Range("O22").Value = "Provincia de Quillota"
$ws.Range("P22").Value = 2700
$ws.Range("Q22").Value = 10
$ws.Range("R22").Value = "Hortaliza"
